$d = $word.ActiveDocument
$CR = [char]13

# ------------------------------------------------------------------
# 1) Insert a new "Meta description" paragraph right after the title
#    (Heading1) paragraph at the top of the document.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titleRange = $titlePara.Range
$titleRange.Collapse(0)           # wdCollapseEnd -> collapse to the end of the title paragraph
$null = $titleRange.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
$metaXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
           '<w:r/>' + `
           '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' + `
           '<w:r><w:t xml:space="preserve">: Explore the Amazon and win big with Amazonia by Merkur. Read our review and play this online slot game for free today.</w:t></w:r>' + `
           '</w:p>'
$null = $metaPara.Range.InsertXML($metaXml)

# ------------------------------------------------------------------
# 2) At the end of the document: remove the duplicated bold title
#    paragraph and rewrite the italic paragraph to hold the image
#    prompt text instead of the meta description text.
# ------------------------------------------------------------------
$boldTitleText = "Play Amazonia Free - Review of Merkur's Online Slot Game"
$italicText = "Explore the Amazon and win big with Amazonia by Merkur. Read our review and play this online slot game for free today."

$boldTitleIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd($CR)
    if ($txt -eq $boldTitleText) {
        $boldTitleIndex = $i
    }
}

if ($boldTitleIndex -gt 0) {
    $d.Paragraphs($boldTitleIndex).Range.Delete()
}

# Re-resolve the italic paragraph index since the document shifted
$italicIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text.TrimEnd($CR)
    if ($txt -eq $italicText) {
        $italicIndex = $i
    }
}

$promptText = "Prompt: Create a cartoon-style feature image for Amazonia slot game. The image should depict a happy Maya warrior with glasses. The background should showcase the lush green of the Amazon rainforest. The Maya warrior should be holding a tablet or smartphone with the Amazonia game logo on it. The image should convey a fun and exciting gaming experience with the Amazonia game in a playful manner. The colours should be bright and vibrant, providing an eye-catching contrast to the green background. Please ensure that the image is in high-resolution to be used not only in the game but for promotional purposes too."

$italicPara = $d.Paragraphs($italicIndex)
$italicRange = $italicPara.Range
$null = $italicRange.MoveEnd(1, -1)
$italicRange.Text = $promptText
